# Requirements List Form - F3: rename sheet, fix print area, adjust view
# zoom and update the footer revision/date stamp (per "feat: sops Update 4").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the main worksheet (S-SW-SC-03 -> F-SW-SD-03). The workbook's
# Print_Area defined name references the sheet by name, so update it too.
$ws.Name = "F-SW-SD-03"
$ws.PageSetup.PrintArea = '$A$1:$D$23'

# Bring the saved view back to a sensible zoom level for the new layout.
$win = $excel.ActiveWindow
$win.Zoom = 80

# Update the footer's revision/date stamp on the right-hand section while
# leaving the left (Issue No.) and center (form code) sections untouched.
$ws.PageSetup.LeftFooter = '&"Arial,Regular"&14Issue No.:(01)'
$ws.PageSetup.CenterFooter = '&"Arial,Regular"&14F-SW-SD/03'
$ws.PageSetup.RightFooter = '&"Arial,Regular"&14Rev:0(01/10/2025)'
